$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last 3 data rows (previously rows 8-10, ECs as sender) entirely
$ws.Range("A8:T10").Delete()

# Update remaining data rows 2-7 with new sender cluster assignments and new TPM-derived values

$ws.Range("A2").Value = "FAPs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.551523333333333
$ws.Range("H2").Value = 4.65457
$ws.Range("I2").Value = 0.950716861801202
$ws.Range("J2").Value = 0.950716861801202
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.4904653333333333
$ws.Range("N2").Value = 1.471396
$ws.Range("O2").Value = 0.1781312055994899
$ws.Range("P2").Value = 0.1781312055994899
$ws.Range("Q2").Value = 0.7609684088577776
$ws.Range("R2").Value = 6.848715679719999
$ws.Range("S2").Value = 0.1693523407764117
$ws.Range("T2").Value = 0.1693523407764117

$ws.Range("A3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.551523333333333
$ws.Range("H3").Value = 4.65457
$ws.Range("I3").Value = 0.950716861801202
$ws.Range("J3").Value = 0.950716861801202
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.266692333333333
$ws.Range("N3").Value = 3.800077
$ws.Range("O3").Value = 0.4600476672363475
$ws.Range("P3").Value = 0.4600476672363475
$ws.Range("Q3").Value = 1.965302711321111
$ws.Range("R3").Value = 17.68772440189
$ws.Range("S3").Value = 0.4373750744739039
$ws.Range("T3").Value = 0.4373750744739039

$ws.Range("A4").Value = "FAPs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.551523333333333
$ws.Range("H4").Value = 4.65457
$ws.Range("I4").Value = 0.950716861801202
$ws.Range("J4").Value = 0.950716861801202
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.996236
$ws.Range("N4").Value = 2.988708
$ws.Range("O4").Value = 0.3618211271641626
$ws.Range("P4").Value = 0.3618211271641626
$ws.Range("Q4").Value = 1.545683399506667
$ws.Range("R4").Value = 13.91115059556
$ws.Range("S4").Value = 0.3439894465508863
$ws.Range("T4").Value = 0.3439894465508864

$ws.Range("A5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.08042766666666666
$ws.Range("H5").Value = 0.241283
$ws.Range("I5").Value = 0.04928313819879805
$ws.Range("J5").Value = 0.04928313819879804
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.4904653333333333
$ws.Range("N5").Value = 1.471396
$ws.Range("O5").Value = 0.1781312055994899
$ws.Range("P5").Value = 0.1781312055994899
$ws.Range("Q5").Value = 0.03944698234088889
$ws.Range("R5").Value = 0.355022841068
$ws.Range("S5").Value = 0.008778864823078169
$ws.Range("T5").Value = 0.008778864823078169

$ws.Range("A6").Value = "MuSCs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.08042766666666666
$ws.Range("H6").Value = 0.241283
$ws.Range("I6").Value = 0.04928313819879805
$ws.Range("J6").Value = 0.04928313819879804
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.266692333333333
$ws.Range("N6").Value = 3.800077
$ws.Range("O6").Value = 0.4600476672363475
$ws.Range("P6").Value = 0.4600476672363475
$ws.Range("Q6").Value = 0.1018771087545556
$ws.Range("R6").Value = 0.916893978791
$ws.Range("S6").Value = 0.02267259276244357
$ws.Range("T6").Value = 0.02267259276244357

$ws.Range("A7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.08042766666666666
$ws.Range("H7").Value = 0.241283
$ws.Range("I7").Value = 0.04928313819879805
$ws.Range("J7").Value = 0.04928313819879804
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.996236
$ws.Range("N7").Value = 2.988708
$ws.Range("O7").Value = 0.3618211271641626
$ws.Range("P7").Value = 0.3618211271641626
$ws.Range("Q7").Value = 0.08012493692933333
$ws.Range("R7").Value = 0.7211244323639999
$ws.Range("S7").Value = 0.01783168061327631
$ws.Range("T7").Value = 0.01783168061327631
